# Auto-generated edit script applying cached-value updates to the Leve profit tables.
# Source: per-sheet H:N (currentAveragePrice.. LeveProfitHQ) recomputation from an external price refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 84.666664
$ws.Range("I9").Value = 80.90000000000001
$ws.Range("J9").Value = 92.2
$ws.Range("K9").Value = 80.90000000000001
$ws.Range("L9").Value = 92.2
$ws.Range("M9").Value = 88.09999999999999
$ws.Range("N9").Value = -430.2
$ws.Range("H98").Value = 3628.9
$ws.Range("J98").Value = 3431.3333
$ws.Range("L98").Value = 3431.3333
$ws.Range("N98").Value = -6427.3333
$ws.Range("H122").Value = 3628.9
$ws.Range("J122").Value = 3431.3333
$ws.Range("L122").Value = 10293.9999
$ws.Range("N122").Value = -15193.9999
$ws.Range("H137").Value = 1863.5238
$ws.Range("I137").Value = 1470.3077
$ws.Range("K137").Value = 4410.9231
$ws.Range("M137").Value = -1860.9231
$ws.Range("H138").Value = 3811.4285
$ws.Range("J138").Value = 3811.4285
$ws.Range("L138").Value = 11434.2855
$ws.Range("N138").Value = -21714.2855

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1883.3334
$ws.Range("I2").Value = 2650
$ws.Range("K2").Value = 2650
$ws.Range("M2").Value = -2537
$ws.Range("H21").Value = 15833.333
$ws.Range("I21").Value = 3750
$ws.Range("J21").Value = 40000
$ws.Range("K21").Value = 3750
$ws.Range("L21").Value = 40000
$ws.Range("M21").Value = -3376
$ws.Range("N21").Value = -40748
$ws.Range("H46").Value = 5027.6
$ws.Range("I46").Value = 5212.6665
$ws.Range("K46").Value = 5212.6665
$ws.Range("M46").Value = -4893.6665
$ws.Range("H61").Value = 2859.5
$ws.Range("I61").Value = 1619
$ws.Range("K61").Value = 1619
$ws.Range("M61").Value = -1407
$ws.Range("H63").Value = 2267.8333
$ws.Range("J63").Value = 3999.5
$ws.Range("L63").Value = 3999.5
$ws.Range("N63").Value = -5371.5
$ws.Range("H66").Value = 2267.8333
$ws.Range("J66").Value = 3999.5
$ws.Range("L66").Value = 19997.5
$ws.Range("N66").Value = -26861.5
$ws.Range("H97").Value = 573.1111
$ws.Range("I97").Value = 452.7143
$ws.Range("J97").Value = 994.5
$ws.Range("K97").Value = 452.7143
$ws.Range("L97").Value = 994.5
$ws.Range("M97").Value = 43.28570000000002
$ws.Range("N97").Value = -1986.5
$ws.Range("H116").Value = 1883.3334
$ws.Range("I116").Value = 2650
$ws.Range("K116").Value = 2650
$ws.Range("M116").Value = -356
$ws.Range("H124").Value = 10707.25
$ws.Range("J124").Value = 10707.25
$ws.Range("L124").Value = 10707.25
$ws.Range("N124").Value = -20527.25
$ws.Range("H136").Value = 2859.5
$ws.Range("I136").Value = 1619
$ws.Range("K136").Value = 4857
$ws.Range("M136").Value = -2307

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1883.3334
$ws.Range("I3").Value = 2650
$ws.Range("K3").Value = 2650
$ws.Range("M3").Value = -2536
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620
$ws.Range("H82").Value = 29123.375
$ws.Range("I82").Value = 10995.667
$ws.Range("J82").Value = 40000
$ws.Range("K82").Value = 10995.667
$ws.Range("L82").Value = 40000
$ws.Range("M82").Value = -10612.667
$ws.Range("N82").Value = -40766
$ws.Range("H85").Value = 29123.375
$ws.Range("I85").Value = 10995.667
$ws.Range("J85").Value = 40000
$ws.Range("K85").Value = 10995.667
$ws.Range("L85").Value = 40000
$ws.Range("M85").Value = -9669.666999999999
$ws.Range("N85").Value = -42652
$ws.Range("H105").Value = 3049
$ws.Range("J105").Value = 3165
$ws.Range("L105").Value = 3165
$ws.Range("N105").Value = -6659
$ws.Range("H107").Value = 1036.625
$ws.Range("I107").Value = 1048.8334
$ws.Range("K107").Value = 1048.8334
$ws.Range("M107").Value = 871.1666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 30832.555
$ws.Range("J59").Value = 34641.43
$ws.Range("L59").Value = 34641.43
$ws.Range("N59").Value = -36931.43
$ws.Range("H60").Value = 17747.666
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26022
$ws.Range("H74").Value = 37880.5
$ws.Range("J74").Value = 37880.5
$ws.Range("L74").Value = 37880.5
$ws.Range("N74").Value = -39628.5
$ws.Range("H77").Value = 37880.5
$ws.Range("J77").Value = 37880.5
$ws.Range("L77").Value = 113641.5
$ws.Range("N77").Value = -122377.5
$ws.Range("H107").Value = 852
$ws.Range("I107").Value = 767.2
$ws.Range("J107").Value = 993.3333
$ws.Range("K107").Value = 767.2
$ws.Range("L107").Value = 993.3333
$ws.Range("M107").Value = 1152.8
$ws.Range("N107").Value = -4833.3333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 98700
$ws.Range("J37").Value = 98700
$ws.Range("L37").Value = 296100
$ws.Range("N37").Value = -296324
$ws.Range("H68").Value = 1816.3334
$ws.Range("I68").Value = 1799
$ws.Range("K68").Value = 5397
$ws.Range("M68").Value = -4586
$ws.Range("H71").Value = 1816.3334
$ws.Range("I71").Value = 1799
$ws.Range("K71").Value = 16191
$ws.Range("M71").Value = -12135
$ws.Range("H81").Value = 6785.9
$ws.Range("I81").Value = 1900
$ws.Range("J81").Value = 8007.375
$ws.Range("K81").Value = 5700
$ws.Range("L81").Value = 24022.125
$ws.Range("M81").Value = -4577
$ws.Range("N81").Value = -26268.125
$ws.Range("H84").Value = 6785.9
$ws.Range("I84").Value = 1900
$ws.Range("J84").Value = 8007.375
$ws.Range("K84").Value = 17100
$ws.Range("L84").Value = 72066.375
$ws.Range("M84").Value = -11484
$ws.Range("N84").Value = -83298.375
$ws.Range("H120").Value = 18333.111
$ws.Range("I120").Value = 10000
$ws.Range("J120").Value = 19374.75
$ws.Range("K120").Value = 30000
$ws.Range("L120").Value = 58124.25
$ws.Range("M120").Value = -25162
$ws.Range("N120").Value = -67800.25
$ws.Range("H131").Value = 2021.8125
$ws.Range("I131").Value = 1279.875
$ws.Range("J131").Value = 2763.75
$ws.Range("K131").Value = 3839.625
$ws.Range("L131").Value = 8291.25
$ws.Range("M131").Value = 1200.375
$ws.Range("N131").Value = -18371.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 38.25
$ws.Range("I2").Value = 38.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 38.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 74.75
$ws.Range("N2").ClearContents()
$ws.Range("H19").Value = 799.5
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H47").Value = 15000
$ws.Range("H80").Value = 3999
$ws.Range("J80").Value = 3999
$ws.Range("L80").Value = 3999
$ws.Range("N80").Value = -5995
$ws.Range("H83").Value = 3999
$ws.Range("J83").Value = 3999
$ws.Range("L83").Value = 19995
$ws.Range("N83").Value = -29979
$ws.Range("H97").Value = 749.94116
$ws.Range("I97").Value = 789.3333
$ws.Range("J97").Value = 454.5
$ws.Range("K97").Value = 789.3333
$ws.Range("L97").Value = 454.5
$ws.Range("M97").Value = -293.3333
$ws.Range("N97").Value = -1446.5
$ws.Range("H126").Value = 1866.6666
$ws.Range("J126").Value = 4000
$ws.Range("L126").Value = 12000
$ws.Range("N126").Value = -16940
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 2288.1428
$ws.Range("I132").Value = 2336.1667
$ws.Range("K132").Value = 7008.500100000001
$ws.Range("M132").Value = -4478.500100000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 701.3333
$ws.Range("I7").Value = 701.3333
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 701.3333
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -589.3333
$ws.Range("N7").ClearContents()
$ws.Range("H54").Value = 48084
$ws.Range("J54").Value = 48084
$ws.Range("L54").Value = 48084
$ws.Range("N54").Value = -49372
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H126").Value = 701.3333
$ws.Range("I126").Value = 701.3333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2103.9999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 366.0001000000002
$ws.Range("N126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 16424.857
$ws.Range("I113").Value = 35001
$ws.Range("J113").Value = 2492.75
$ws.Range("K113").Value = 105003
$ws.Range("L113").Value = 7478.25
$ws.Range("M113").Value = -102833
$ws.Range("N113").Value = -11818.25
$ws.Range("H122").Value = 837.25
$ws.Range("I122").Value = 837.25
$ws.Range("K122").Value = 2511.75
$ws.Range("M122").Value = -61.75
$ws.Range("H126").Value = 4288.5
$ws.Range("I126").Value = 4288.5
$ws.Range("K126").Value = 12865.5
$ws.Range("M126").Value = -10395.5
